$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "2022-Q3" worksheet.
#    Duplicate the "2021-Q4" sheet (so formatting/styles/column layout match
#    the other quarterly sheets exactly) and drop the copy in front of it.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($q4)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Fill in the new quarter's fund holdings. Columns B-G hold numeric-looking
# strings that must stay text (as in every other quarterly sheet), so force
# a text number format before writing them, then restore the default style.
$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("B2").Value = "014126"
$q3.Range("C2").Value = "华夏中证1000指数增强C"
$q3.Range("D2").Value = "8.78"
$q3.Range("E2").Value = "89.62"
$q3.Range("F2").Value = "0.82"
$q3.Range("G2").Value = "0.0720"

$q3.Range("B3").Value = "014125"
$q3.Range("C3").Value = "华夏中证1000指数增强A"
$q3.Range("D3").Value = "0.97"
$q3.Range("E3").Value = "89.62"
$q3.Range("F3").Value = "0.82"
$q3.Range("G3").Value = "0.0080"

$q3.Range("B2:G3").Style = "Normal"

# Column H ("仓位排名") is numeric - write it after restoring the style so it
# isn't coerced into a text cell.
$q3.Range("H2").Value = 4
$q3.Range("H3").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the new 2022-Q3 row at the
#    top of the data and push the existing rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for the extra row: copy row 5's style onto the new row 6 first.
$total.Range("A5").Copy($total.Range("A6"))

$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 6
$total.Range("D6").Value = 1.79

$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 1.27

$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 1.07

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.09

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.08

$total.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 3. Restore the originally active tab ("2020-Q4"); copying a sheet makes
#    the new copy active, which we don't want to leave as a side effect.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()

